# "contingencies with rene fine"
# Populate B1 / A2 with a 0 value (styled: bold, centered/top, thin black box
# border) and B2 with the label "disconnected_elements".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Build the styled format once on B1 ...
$r1 = $ws.Range("B1")
$r1.Borders.Color = 0
$r1.Borders.LineStyle = 1
$r1.Borders.Weight = 2
$r1.Font.Bold = $true
$r1.HorizontalAlignment = -4108
$r1.VerticalAlignment = -4160

# ... then copy the same format onto A2 so both cells share one style
# definition instead of each independently rebuilding the style table.
$r1.Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
